$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert 7 new rows before the existing row 4 (logF_BusinessProcessName row),
# pushing it (and everything below) down to row 11.
$ws.Rows("4:10").Insert()

# --- Email Credentials block (rows 6-10), filled column by column ---
$ws.Range("A6").Value = "EmailCredentialAssetName"
$ws.Range("A7").Value = "EmailCredentialFolderPath"
$ws.Range("A8").Value = "EmailServer"
$ws.Range("A9").Value = "EmailPort"
$ws.Range("A10").Value = "OverrideEmail"

$ws.Range("B6").Value = "Email Credentials"
$ws.Range("B7").Value = "P3 Automation"
$ws.Range("B8").Value = "smtp.gmail.com"
$ws.Range("B9").Value = 465
$ws.Range("B10").Value = "wplee.327@gmail.com"

$ws.Range("C6").Value = "Name of Orchestrator asset containing email credentials."
$ws.Range("C7").Value = "Path to Orchestrator folder containing email credentials."
$ws.Range("C8").Value = "Server for email origin."
$ws.Range("C9").Value = "Port for email origin."
$ws.Range("C10").Value = "Override email(s) that reports are sent to. Comma separated."

# --- GitHub Credentials block (rows 4-5) ---
$ws.Range("A4").Value = "GitHubCredentialAssetName"
$ws.Range("A5").Value = "GitHubCredentialFolderPath"

$ws.Range("C4").Value = "Name of Orchestrator asset containing GitHub credentials."
$ws.Range("C5").Value = "Path to Orchestrator folder containing GitHub credentials."

$ws.Range("B4").Value = "Git Credentials"
$ws.Range("B5").Value = "P3 Automation"

# Restore the default custom row height on the newly inserted rows.
$ws.Rows("4:10").RowHeight = 14.25

# Update the selected cell on the worksheet.
$ws.Range("B5").Select()

Write-Output "done"
